# Fruta / hortaliza, semanal
# Apply the weekly refresh of the "Higo" (fig) sheet: several report rows
# get reshuffled to different dates (the underlying weekly data block for
# each date moved to a different row position), with their Volumen,
# Precio minimo/maximo/promedio, Origen and Precio $/Kg following along.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Primera, was 2021-04-... serial 44301) -> now date 44302
$ws.Cells.Item(2, 4).Value = 44302
$ws.Cells.Item(2, 13).Value = 50
$ws.Cells.Item(2, 14).Value = 15000
$ws.Cells.Item(2, 15).Value = 15000
$ws.Cells.Item(2, 16).Value = 15000
$ws.Cells.Item(2, 19).Value = 2143

# Row 3 (Segunda, paired with row 2)
$ws.Cells.Item(3, 4).Value = 44302
$ws.Cells.Item(3, 13).Value = 30

# Row 6 (Primera) -> now date 44292
$ws.Cells.Item(6, 4).Value = 44292
$ws.Cells.Item(6, 13).Value = 25
$ws.Cells.Item(6, 14).Value = 16000
$ws.Cells.Item(6, 15).Value = 16000
$ws.Cells.Item(6, 16).Value = 16000
$ws.Cells.Item(6, 19).Value = 2286

# Row 7 (Segunda, paired with row 6)
$ws.Cells.Item(7, 4).Value = 44292
$ws.Cells.Item(7, 14).Value = 15000
$ws.Cells.Item(7, 15).Value = 15000
$ws.Cells.Item(7, 16).Value = 15000
$ws.Cells.Item(7, 19).Value = 2143

# Row 8 (Primera) -> now date 44299, Origen changes to Provincia de Santiago
$ws.Cells.Item(8, 4).Value = 44299
$ws.Cells.Item(8, 13).Value = 80
$ws.Cells.Item(8, 14).Value = 15000
$ws.Cells.Item(8, 15).Value = 15000
$ws.Cells.Item(8, 16).Value = 15000
$ws.Cells.Item(8, 18).Value = "Provincia de Santiago"
$ws.Cells.Item(8, 19).Value = 2143

# Row 9 (Segunda, paired with row 8)
$ws.Cells.Item(9, 4).Value = 44299
$ws.Cells.Item(9, 13).Value = 75
$ws.Cells.Item(9, 14).Value = 12000
$ws.Cells.Item(9, 15).Value = 12000
$ws.Cells.Item(9, 16).Value = 12000
$ws.Cells.Item(9, 18).Value = "Provincia de Santiago"
$ws.Cells.Item(9, 19).Value = 1714

# Row 10 (Primera) -> now date 44320, Origen changes to Región Metropolitana
$ws.Cells.Item(10, 4).Value = 44320
$ws.Cells.Item(10, 13).Value = 20
$ws.Cells.Item(10, 14).Value = 12000
$ws.Cells.Item(10, 15).Value = 12000
$ws.Cells.Item(10, 16).Value = 12000
$ws.Cells.Item(10, 18).Value = "Región Metropolitana"
$ws.Cells.Item(10, 19).Value = 1714

# Row 11 (Segunda, paired with row 10)
$ws.Cells.Item(11, 4).Value = 44320
$ws.Cells.Item(11, 13).Value = 30
$ws.Cells.Item(11, 14).Value = 8000
$ws.Cells.Item(11, 15).Value = 8000
$ws.Cells.Item(11, 16).Value = 8000
$ws.Cells.Item(11, 18).Value = "Región Metropolitana"
$ws.Cells.Item(11, 19).Value = 1143

# Row 12 (Primera) -> now date 44301
$ws.Cells.Item(12, 4).Value = 44301
$ws.Cells.Item(12, 13).Value = 100
$ws.Cells.Item(12, 14).Value = 14000
$ws.Cells.Item(12, 15).Value = 14000
$ws.Cells.Item(12, 16).Value = 14000
$ws.Cells.Item(12, 19).Value = 2000

# Row 13 (Segunda, paired with row 12)
$ws.Cells.Item(13, 4).Value = 44301
$ws.Cells.Item(13, 14).Value = 12000
$ws.Cells.Item(13, 15).Value = 12000
$ws.Cells.Item(13, 16).Value = 12000
$ws.Cells.Item(13, 19).Value = 1714

# Row 14 (Primera) -> now date 44322
$ws.Cells.Item(14, 4).Value = 44322
$ws.Cells.Item(14, 13).Value = 45
$ws.Cells.Item(14, 14).Value = 12000
$ws.Cells.Item(14, 15).Value = 12000
$ws.Cells.Item(14, 16).Value = 12000
$ws.Cells.Item(14, 19).Value = 1714

# Row 15 (Segunda, paired with row 14)
$ws.Cells.Item(15, 4).Value = 44322
$ws.Cells.Item(15, 13).Value = 80
$ws.Cells.Item(15, 14).Value = 8000
$ws.Cells.Item(15, 15).Value = 8000
$ws.Cells.Item(15, 16).Value = 8000
$ws.Cells.Item(15, 19).Value = 1143
